# Update Work Week and Social Spending
# (commit message refers to upstream source data refresh; the actual edit
#  updates the Lithuania GDP per Capita "Data" sheet with revised figures
#  and extends the series through 2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Map of row number -> new "Data" column (E) value. These are stored as
# text in the source workbook (shared strings), not numbers, even though
# most look numeric - so NumberFormat is forced to Text ("@") on each cell
# before the value is written, which is what keeps Excel from silently
# re-typing the cell as a number.
$dataValues = [ordered]@{
    2  = "12103"
    9  = "12734"
    10 = "13190"
    11 = "13654"
    12 = "13718"
    13 = "13796"
    14 = "13287"
    15 = "13667"
    16 = "13831"
    17 = "14851"
    18 = "14693"
    19 = "13809"
    20 = "12999.308856948"
    21 = "10240.9392301755"
    22 = "8621.09647172474"
    23 = "7830.00513468936"
    24 = "8299.46993245928"
    25 = "8793.70079785355"
    26 = "9592.63965461508"
    27 = "10382.117408001"
    28 = "10336.3951624259"
    29 = "10806.8738960776"
    30 = "11606.728929678"
    31 = "12490.6211630176"
    32 = "13917.5076600235"
    33 = "14995.7669761907"
    34 = "16417.9484262796"
    35 = "17916.0591594461"
    36 = "20138.4724101225"
    37 = "20879.8024701771"
    38 = "17983.3490419282"
    39 = "18663.7609461325"
    40 = "20243"
    41 = "21303"
    42 = "22274"
    43 = "23252"
    44 = "23889"
    45 = "24743"
}

# New rows (2011-2016) need the rest of their columns populated too, same
# as every other existing row: Country Code / Country Name / Indicator / Year.
$newRows = [ordered]@{
    40 = 2011
    41 = 2012
    42 = 2013
    43 = 2014
    44 = 2015
    45 = 2016
}

foreach ($row in $newRows.Keys) {
    $ws.Range("A$row").Value = 440
    $ws.Range("B$row").Value = "Lithuania"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $newRows[$row]
}

foreach ($row in $dataValues.Keys) {
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dataValues[$row]
}
